$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E width change (20.28515625 -> 21) ---
$ws.Columns("E").ColumnWidth = 21

# --- Row 35: new player-feedback note + assignee, with taller row ---
$ws.Rows(35).RowHeight = 60
$ws.Range("E35").Value = "player feedback: impacts, camera shake - do through sounds & vfx"
$ws.Range("F11").Copy()
$ws.Range("F35").PasteSpecial(-4122)
$ws.Range("F35").Value = "Charlie"

# --- Row 37: "Testing" row now assigned to everyone ---
$ws.Range("F3").Copy()
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("F37").Value = "everyone"

# --- Row 38: feedback note rewording ---
$ws.Range("E38").Value = "Iterate based on testing feedback"

# --- Update selection to match author's final cursor position ---
$null = $ws.Range("D38").Select()

$excel.CutCopyMode = $false
